# Update the division problems in the practice table.
# Each (old, new) pair below is unique in the document, so a single
# Content-wide Find/ReplaceAll per pair is safe. The one pair whose
# *result* text ("10÷8=") collides with another pair's *source* text
# ("87÷9=" -> "10÷8=") is ordered so the source->10÷8= rewrite happens
# before the newly-created "10÷8=" cell is produced, avoiding a
# double-replace.
$d = $word.ActiveDocument

$replacements = @(
    @("61÷4=", "50÷8="),
    @("68÷6=", "48÷8="),
    @("17÷6=", "88÷8="),
    @("94÷9=", "94÷3="),
    @("94÷2=", "61÷9="),
    @("57÷7=", "77÷4="),
    @("10÷8=", "37÷4="),
    @("87÷9=", "10÷8="),
    @("31÷6=", "32÷8="),
    @("66÷9=", "17÷2="),
    @("93÷6=", "12÷9="),
    @("72÷6=", "25÷8="),
    @("64÷7=", "62÷9="),
    @("95÷4=", "66÷3="),
    @("39÷2=", "77÷9="),
    @("69÷7=", "95÷8="),
    @("42÷6=", "10÷3="),
    @("65÷7=", "75÷2="),
    @("98÷5=", "60÷5="),
    @("59÷7=", "28÷9="),
    @("79÷5=", "14÷3="),
    @("60÷6=", "79÷3="),
    @("53÷2=", "26÷5="),
    @("79÷2=", "47÷3="),
    @("91÷3=", "19÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
